$wb = $excel.ActiveWorkbook

# --- Update "Alternative_Input" sheet (data + selection) ---
$wsAlt = $wb.Worksheets.Item("Alternative_Input")

# New data for A1:C7 (copied from "FImW single Jenal Fluo")
$wsAlt.Range("A1").Value = "1634 mNG_FimW"
$wsAlt.Range("B1").Value = 20220726
$wsAlt.Range("C1").Value = "5s interval-2h37"

$wsAlt.Range("A2").Value = "1634 mNG_FimW"
$wsAlt.Range("B2").Value = 20220804
$wsAlt.Range("C2").Value = "5s interval-2h37"

$wsAlt.Range("A3").Value = "1635 mNG_FimX"
$wsAlt.Range("B3").Value = 20220726
$wsAlt.Range("C3").Value = "5s interval-2h37"

$wsAlt.Range("A4").Value = "1635 mNG_FimX"
$wsAlt.Range("B4").Value = 20220804
$wsAlt.Range("C4").Value = "5s interval-2h37"

$wsAlt.Range("A5").Value = "1638 mNG_FimW pch-"
$wsAlt.Range("B5").Value = 20220726
$wsAlt.Range("C5").Value = "5s interval-2h37"

$wsAlt.Range("A6").Value = "1638 mNG_FimW pch-"
$wsAlt.Range("B6").Value = 20220728
$wsAlt.Range("C6").Value = "5s interval-2h37"

$wsAlt.Range("A7").Value = "1638 mNG_FimW pch-"
$wsAlt.Range("B7").Value = 20220729
$wsAlt.Range("C7").Value = "5s interval-2h37"

# Clear old rows 8:10 (keep formatting, remove values)
$wsAlt.Range("A8:C10").ClearContents()

# Touch formatting (no-op) on the non-styled A/C cells so the engine keeps
# an explicit empty <c> element for them in the saved XML (matches target).
$wsAlt.Range("A8").Font.Bold = $false
$wsAlt.Range("C8").Font.Bold = $false
$wsAlt.Range("A9").Font.Bold = $false
$wsAlt.Range("C9").Font.Bold = $false
$wsAlt.Range("A10").Font.Bold = $false
$wsAlt.Range("C10").Font.Bold = $false

# Update selection on this sheet
$wsAlt.Range("C16").Select()

# --- Update "FImW single Jenal Fluo" sheet selection only ---
$wsFimW = $wb.Worksheets.Item("FImW single Jenal Fluo")
$wsFimW.Activate()
$wsFimW.Range("A1:C7").Select()

$wsAlt.Activate()
